$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new price text is a "pure number" token (e.g. "308.07")
# must have their NumberFormat switched to Text first, otherwise Excel COM
# auto-converts the assigned string into a floating point number (losing
# trailing/leading zeros) exactly like typing it into a General-formatted cell.
$textPriceRows = @(6,7,8,9,10,11,12,13,15,16,17,18,19,21,23,24,27,28,29,30,31,32,33,34,36,38,39,40,41,42,43,45,46,47,48,49,50,51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.408.82"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "1.796.29"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "308.07"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "0.4530"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("D8").Value = "0.3594"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "46.35"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").Value = "0.07122"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "0.8873"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "0.07829"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "19.52"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "1.833.71"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "5.282"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "6.336"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "84.98"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "0.000008583"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "14.30"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Value = "26.445.41"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").Value = "4.992"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "10.53"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.018.05"
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "152.64"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").Value = "17.94"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").Value = "2.049"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("D30").Value = "112.10"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "4.881"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "0.08662"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "3.048"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").Value = "2.757"
$ws.Range("E34").Value = "  +6.84%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "0.7287"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "1.076"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "0.01941"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").Value = "0.05111"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "2.878"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "0.5173"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("D43").Value = "6.909"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("D45").Value = "8.015"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("D46").Value = "0.4675"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "9.877"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "100.67"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "1.587"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "64.53"
$ws.Range("E51").Value = "  +0.22%  "
